# Add two new academic-group rows (Магистратура / Управление в технических
# системах / Беспилотная робототехника) to "Лист1": 201-321 and 201-322,
# both sharing the same academic plan text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# Use the last existing data row (38) as the formatting template for the
# two new rows, then overwrite the values that differ.
$template = $ws.Range("A38:K38")

$row39 = $ws.Range("A39:K39")
$template.Copy($row39)

$row40 = $ws.Range("A40:K40")
$template.Copy($row40)

# Row 39: group 201-321
$ws.Range("A39").Value = 38
$ws.Range("B39").Value = "201-321"
$ws.Range("C39").Value = "Академический учебный план 000019339 от 31.01.2024 11:51:06"
$ws.Range("D39").Value = "Факультет информационных технологий"
$ws.Range("E39").Value = "27.04.04"
$ws.Range("F39").Value = "Управление в технических системах"
$ws.Range("G39").Value = "Очная"
$ws.Range("H39").Value = "Магистратура"
$ws.Range("I39").Value = "Беспилотная робототехника"
$ws.Range("J39").Value = "Первый"
$ws.Range("K39").Value = "Нет"

# Row 40: group 201-322
$ws.Range("A40").Value = 39
$ws.Range("B40").Value = "201-322"
$ws.Range("C40").Value = "Академический учебный план 000019339 от 31.01.2024 11:51:06"
$ws.Range("D40").Value = "Факультет информационных технологий"
$ws.Range("E40").Value = "27.04.04"
$ws.Range("F40").Value = "Управление в технических системах"
$ws.Range("G40").Value = "Очная"
$ws.Range("H40").Value = "Магистратура"
$ws.Range("I40").Value = "Беспилотная робототехника"
$ws.Range("J40").Value = "Первый"
$ws.Range("K40").Value = "Нет"
